$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    ,@("Login with valid username and password", "PASSED", "chrome")
    ,@("Fees create and delete functionality", "FAILED", "chrome")
    ,@("Create a country", "FAILED", "chrome")
    ,@("Login with valid username and password", "PASSED", "chrome")
    ,@("Create a country", "PASSED", "chrome")
    ,@("Fees create and delete functionality", "FAILED", "chrome")
    ,@("Fees create and delete functionality", "PASSED", "chrome")
    ,@("Login with valid username and password", "PASSED", "chrome")
    ,@("Create a country", "PASSED", "chrome")
    ,@("Fees create and delete functionality", "PASSED", "chrome")
    ,@("Fees create and delete functionality", "PASSED", "chrome")
    ,@("Login with valid username and password", "PASSED", "chrome")
    ,@("Login with valid username and password", "PASSED", "chrome")
    ,@("Create a country", "PASSED", "chrome")
    ,@("Fees create and delete functionality", "PASSED", "chrome")
)

$startRow = 551
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
